$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3561.375
$ws.Range("J112").Value = 3562.4092
$ws.Range("L112").Value = 10687.2276
$ws.Range("N112").Value = -12903.2276
$ws.Range("H125").Value = 5450.2
$ws.Range("I125").Value = 2252.4443
$ws.Range("K125").Value = 20271.9987
$ws.Range("M125").Value = -17811.9987
$ws.Range("H129").Value = 1847.1632
$ws.Range("I129").Value = 1528.6666
$ws.Range("K129").Value = 4585.9998
$ws.Range("M129").Value = 414.0002000000004
$ws.Range("H137").Value = 12263905
$ws.Range("I137").Value = 1112240.4
$ws.Range("K137").Value = 3336721.2
$ws.Range("M137").Value = -3334171.2
$ws.Range("H138").Value = 1827.2526
$ws.Range("I138").Value = 969.9286
$ws.Range("J138").Value = 1968.4589
$ws.Range("K138").Value = 2909.7858
$ws.Range("L138").Value = 5905.376700000001
$ws.Range("M138").Value = 2230.2142
$ws.Range("N138").Value = -16185.3767

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3434.14
$ws.Range("I32").Value = 3189.1685
$ws.Range("J32").Value = 8088.6
$ws.Range("K32").Value = 3189.1685
$ws.Range("L32").Value = 8088.6
$ws.Range("M32").Value = -2902.1685
$ws.Range("N32").Value = -8662.6
$ws.Range("H45").Value = 3235
$ws.Range("I45").Value = 3124.4285
$ws.Range("J45").Value = 3345.5715
$ws.Range("K45").Value = 3124.4285
$ws.Range("L45").Value = 3345.5715
$ws.Range("M45").Value = -2747.4285
$ws.Range("N45").Value = -4099.5715
$ws.Range("H61").Value = 3586.5588
$ws.Range("I61").Value = 2845.5454
$ws.Range("K61").Value = 2845.5454
$ws.Range("M61").Value = -2633.5454
$ws.Range("H74").Value = 11906830
$ws.Range("J74").Value = 3624.5
$ws.Range("L74").Value = 3624.5
$ws.Range("N74").Value = -5372.5
$ws.Range("H77").Value = 11906830
$ws.Range("J77").Value = 3624.5
$ws.Range("L77").Value = 18122.5
$ws.Range("N77").Value = -26858.5
$ws.Range("H132").Value = 3456.45
$ws.Range("I132").Value = 1794.2593
$ws.Range("K132").Value = 5382.7779
$ws.Range("M132").Value = -2852.7779
$ws.Range("H136").Value = 3586.5588
$ws.Range("I136").Value = 2845.5454
$ws.Range("K136").Value = 8536.636200000001
$ws.Range("M136").Value = -5986.636200000001
$ws.Range("H139").Value = 110500
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 110500
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 110500
$ws.Range("N139").Value = -120780
$ws.Range("M139").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 359
$ws.Range("I22").Value = 245.5
$ws.Range("K22").Value = 245.5
$ws.Range("M22").Value = -72.5
$ws.Range("H140").Value = 170143.42
$ws.Range("J140").Value = 170143.42
$ws.Range("L140").Value = 170143.42
$ws.Range("N140").Value = -180503.42

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7460.7856
$ws.Range("J99").Value = 4469.5293
$ws.Range("L99").Value = 4469.5293
$ws.Range("N99").Value = -7465.5293
$ws.Range("H126").Value = 7460.7856
$ws.Range("J126").Value = 4469.5293
$ws.Range("L126").Value = 13408.5879
$ws.Range("N126").Value = -18348.5879
$ws.Range("H134").Value = 2270.2122
$ws.Range("I134").Value = 2280.742
$ws.Range("K134").Value = 6842.226000000001
$ws.Range("M134").Value = -4307.226000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1211.3334
$ws.Range("I5").Value = 546.3333
$ws.Range("K5").Value = 1638.9999
$ws.Range("M5").Value = -1526.9999
$ws.Range("H135").Value = 1211.3334
$ws.Range("I135").Value = 546.3333
$ws.Range("K135").Value = 4916.9997
$ws.Range("M135").Value = -2381.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5584.3
$ws.Range("I70").Value = 5000.636
$ws.Range("K70").Value = 5000.636
$ws.Range("M70").Value = -4730.636
$ws.Range("H73").Value = 5584.3
$ws.Range("I73").Value = 5000.636
$ws.Range("K73").Value = 5000.636
$ws.Range("M73").Value = -4064.636
$ws.Range("H80").Value = 88452.664
$ws.Range("I80").Value = 129948.375
$ws.Range("J80").Value = 5461.25
$ws.Range("K80").Value = 129948.375
$ws.Range("L80").Value = 5461.25
$ws.Range("M80").Value = -128950.375
$ws.Range("N80").Value = -7457.25
$ws.Range("H83").Value = 88452.664
$ws.Range("I83").Value = 129948.375
$ws.Range("J83").Value = 5461.25
$ws.Range("K83").Value = 649741.875
$ws.Range("L83").Value = 27306.25
$ws.Range("M83").Value = -644749.875
$ws.Range("N83").Value = -37290.25
$ws.Range("H123").Value = 49296.152
$ws.Range("J123").Value = 49296.152
$ws.Range("L123").Value = 49296.152
$ws.Range("N123").Value = -54196.152
$ws.Range("H126").Value = 2743.9375
$ws.Range("I126").Value = 2149.1292
$ws.Range("J126").Value = 3828.5881
$ws.Range("K126").Value = 6447.3876
$ws.Range("L126").Value = 11485.7643
$ws.Range("M126").Value = -3977.3876
$ws.Range("N126").Value = -16425.7643
$ws.Range("H132").Value = 107600.945
$ws.Range("I132").Value = 201148.5
$ws.Range("K132").Value = 603445.5
$ws.Range("M132").Value = -600915.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1114
$ws.Range("J22").Value = 1333.1666
$ws.Range("L22").Value = 1333.1666
$ws.Range("N22").Value = -1923.1666
$ws.Range("H27").Value = 1114
$ws.Range("J27").Value = 1333.1666
$ws.Range("L27").Value = 1333.1666
$ws.Range("N27").Value = -1547.1666
$ws.Range("H68").Value = 6945.0586
$ws.Range("J68").Value = 8275.857
$ws.Range("L68").Value = 8275.857
$ws.Range("N68").Value = -9773.857
$ws.Range("H71").Value = 6945.0586
$ws.Range("J71").Value = 8275.857
$ws.Range("L71").Value = 41379.285
$ws.Range("N71").Value = -48867.285
$ws.Range("H100").Value = 3489.875
$ws.Range("I100").Value = 3184.2
$ws.Range("J100").Value = 3999.3333
$ws.Range("K100").Value = 3184.2
$ws.Range("L100").Value = 3999.3333
$ws.Range("M100").Value = -2643.2
$ws.Range("N100").Value = -5081.3333
$ws.Range("H122").Value = 9118.947
$ws.Range("I122").Value = 1718.1666
$ws.Range("J122").Value = 12534.692
$ws.Range("K122").Value = 5154.4998
$ws.Range("L122").Value = 37604.076
$ws.Range("M122").Value = -2704.4998
$ws.Range("N122").Value = -42504.076

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1547.4615
$ws.Range("J96").Value = 2324.6667
$ws.Range("L96").Value = 2324.6667
$ws.Range("N96").Value = -5070.6667
$ws.Range("H100").Value = 10899.728
$ws.Range("I100").Value = 1975
$ws.Range("J100").Value = 15999.571
$ws.Range("K100").Value = 3950
$ws.Range("L100").Value = 31999.142
$ws.Range("M100").Value = -3409
$ws.Range("N100").Value = -33081.142
$ws.Range("H122").Value = 2999.9524
$ws.Range("I122").Value = 2824
$ws.Range("J122").Value = 3351.8572
$ws.Range("K122").Value = 8472
$ws.Range("L122").Value = 10055.5716
$ws.Range("M122").Value = -6022
$ws.Range("N122").Value = -14955.5716
$ws.Range("H136").Value = 3247.726
$ws.Range("I136").Value = 2113.9648
$ws.Range("J136").Value = 5641.222
$ws.Range("K136").Value = 6341.894400000001
$ws.Range("L136").Value = 16923.666
$ws.Range("M136").Value = -3791.894400000001
$ws.Range("N136").Value = -22023.666
